# Error Calculations and Plots
# Apply edits to missing_data.xlsx worksheet:
#  1. Remove two rows ("RM 232" at row 26 and, after the shift, "SC 92")
#     which shifts all subsequent rows up by two.
#  2. Clear / fill scattered values in column F.
#  3. Clear / fill a few remaining B / F values that differ after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the two rows that were removed from the data set ---
# Row 26 is "RM 232" in the original layout.
$ws.Rows("26:26").Delete()
# After that delete, the row that was "SC 92" (originally row 28) is now row 27.
$ws.Rows("27:27").Delete()

# --- Step 2: column F value swaps (values newly imputed / newly cleared) ---
$ws.Range("F6").Value2 = 16.43
$ws.Range("F8").ClearContents()
$ws.Range("F12").Value2 = 17.45
$ws.Range("F14").ClearContents()
$ws.Range("F17").Value2 = 17.78
$ws.Range("F18").Value2 = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("F23").Value2 = 16.48

# --- Step 3: remaining cell-level corrections on the shifted rows ---
# Row 27 = "SC 101"
$ws.Range("B27").Value2 = -20.4
$ws.Range("F27").ClearContents()
# Row 29 = "SC 119"
$ws.Range("B29").ClearContents()
# Row 32 = "SC 193"
$ws.Range("B32").ClearContents()
